# Add response and question about place to pay taxes
#
# 1) Insert 3 new "greeting" rows (responses) right after the existing
#    greeting block (old row 8), before the "invite-eating" block that used
#    to start at row 9.
# 2) Insert 2 new "law-place" rows (questions) at the end of the existing
#    "law-place" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: three new greeting rows, inserted before old row 9 ---
$ws.Rows("9:11").Insert()

$ws.Cells.Item(9, 1).Value = "greeting"
$ws.Cells.Item(9, 2).Value = "หวัดดีค่ะ"

$ws.Cells.Item(10, 1).Value = "greeting"
$ws.Cells.Item(10, 2).Value = "หวัดดีคับ"

$ws.Cells.Item(11, 1).Value = "greeting"
$ws.Cells.Item(11, 2).Value = "หวัดดีครับ"

# --- Step 2: two new law-place rows, appended after the old last
#     law-place row (row 90), which is now row 93 after the insert above ---
$ws.Rows("94:95").Insert()

$ws.Cells.Item(94, 1).Value = "law-place"
$ws.Cells.Item(94, 2).Value = "เสียภาษีไหน"

$ws.Cells.Item(95, 1).Value = "law-place"
$ws.Cells.Item(95, 2).Value = "จ่ายภาษีไหน"

# --- Restore the view state captured in the edited workbook ---
$ws.Range("I63").Select()
